$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A14").Value = 44616
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat
$ws.Range("B14").Value = 1.25
$ws.Range("B14").HorizontalAlignment = $ws.Range("B13").HorizontalAlignment
$ws.Range("B14").VerticalAlignment = $ws.Range("B13").VerticalAlignment
$ws.Range("C14").Value = "Adding another scatter plot, highlighting outliers"

$ws.Range("C14").Select()
